$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '30.298.50'
Set-TextValue 'E2' '  +0.21%  '
Set-TextValue 'D3' '1.869.41'
Set-TextValue 'E3' '  +0.39%  '
Set-TextValue 'E4' '  -0.14%  '
Set-TextValue 'D5' '235.07'
Set-TextValue 'E5' '  -0.63%  '
Set-TextValue 'E6' '  -0.05%  '
Set-TextValue 'D7' '0.4695'
Set-TextValue 'E7' '  +0.42%  '
Set-TextValue 'D8' '0.2869'
Set-TextValue 'E8' '  +0.50%  '
Set-TextValue 'D9' '0.06595'
Set-TextValue 'E9' '  +1.03%  '
Set-TextValue 'D10' '21.74'
Set-TextValue 'E10' '  -1.10%  '
Set-TextValue 'D11' '0.07964'
Set-TextValue 'E11' '  +0.62%  '
Set-TextValue 'D12' '96.79'
Set-TextValue 'E12' '  -0.71%  '
Set-TextValue 'D13' '1.875.69'
Set-TextValue 'E13' '  +0.69%  '
Set-TextValue 'D14' '0.6988'
Set-TextValue 'E14' '  +2.58%  '
Set-TextValue 'D15' '5.116'
Set-TextValue 'E15' '  -0.87%  '
Set-TextValue 'D16' '268.56'
Set-TextValue 'E16' '  -1.03%  '
Set-TextValue 'D17' '30.343.38'
Set-TextValue 'E17' '  +0.36%  '
Set-TextValue 'D18' '14.15'
Set-TextValue 'E18' '  +4.59%  '
Set-TextValue 'D19' '0.000007787'
Set-TextValue 'E19' '  +6.02%  '
Set-TextValue 'E20' '  -0.01%  '
Set-TextValue 'E21' '  +0.42%  '
Set-TextValue 'E22' '  -0.24%  '
Set-TextValue 'D23' '5.265'
Set-TextValue 'E23' '  -1.06%  '
Set-TextValue 'D24' '6.219'
Set-TextValue 'E24' '  +0.68%  '
Set-TextValue 'D25' '9.391'
Set-TextValue 'E25' '  +1.88%  '
Set-TextValue 'D26' '167.43'
Set-TextValue 'E26' '  -0.21%  '
Set-TextValue 'D27' '18.87'
Set-TextValue 'E27' '  -0.37%  '
Set-TextValue 'D28' '1.953'
Set-TextValue 'E28' '  +0.23%  '
Set-TextValue 'E29' '  -1.52%  '
Set-TextValue 'D30' '0.09889'
Set-TextValue 'E30' '  +0.61%  '
Set-TextValue 'D31' '4.347'
Set-TextValue 'E31' '  -0.27%  '
Set-TextValue 'D32' '1.458'
Set-TextValue 'E32' '  -1.65%  '
Set-TextValue 'D33' '4.056'
Set-TextValue 'E33' '  +0.00%  '
Set-TextValue 'D34' '0.04718'
Set-TextValue 'E34' '  +0.15%  '
Set-TextValue 'D35' '1.137'
Set-TextValue 'E35' '  +0.39%  '
Set-TextValue 'D36' '0.7028'
Set-TextValue 'E36' '  +0.16%  '
Set-TextValue 'D37' '2.724'
Set-TextValue 'E37' '  +0.53%  '
Set-TextValue 'D38' '0.01876'
Set-TextValue 'E38' '  +0.09%  '
Set-TextValue 'D39' '2.807'
Set-TextValue 'E39' '  +6.78%  '
Set-TextValue 'D40' '6.261'
Set-TextValue 'E40' '  -0.20%  '
Set-TextValue 'D41' '72.30'
Set-TextValue 'E41' '  -4.17%  '
Set-TextValue 'D42' '1.960'
Set-TextValue 'E42' '  +0.83%  '
Set-TextValue 'D43' '0.8443'
Set-TextValue 'E43' '  -0.85%  '
Set-TextValue 'D44' '0.4177'
Set-TextValue 'E44' '  +0.46%  '
Set-TextValue 'E45' '  -0.01%  '
Set-TextValue 'D46' '102.91'
Set-TextValue 'E46' '  -0.17%  '
Set-TextValue 'D47' '7.106'
Set-TextValue 'E47' '  -0.94%  '
Set-TextValue 'D48' '9.115'
Set-TextValue 'E48' '  -1.66%  '
Set-TextValue 'D49' '920.98'
Set-TextValue 'E49' '  -3.16%  '
Set-TextValue 'D50' '34.63'
Set-TextValue 'E50' '  +1.49%  '
Set-TextValue 'D51' '0.05689'
Set-TextValue 'E51' '  +0.67%  '

Write-Output "Applied cryptos update"